$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) — update F3:F6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 310
$ws1.Range("F4").Value = 1312
$ws1.Range("F5").Value = 85
$ws1.Range("F6").Value = 65

# Sheet "全部类型" (all types) — update F4:F7
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 310
$ws4.Range("F5").Value = 1312
$ws4.Range("F6").Value = 85
$ws4.Range("F7").Value = 65
